# Add data for 2022-07-23
# Updates the "through" date in the sheet name and header string, and
# updates/adds carjacking counts for column B (July 2022 through <date>)
# plus a handful of historical corrections scattered across the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update the header label for the current month column.
$ws.Name = "Through 2022-07-15"
$ws.Range("B1").Value = "July 2022 (through July 15)"

# Cell updates: Row => { Column => NewValue }
$updates = @{
    2  = @{ B = 9; I = 8; AD = 7; AK = 3; AY = 3 }
    3  = @{ AR = 5 }
    4  = @{ B = 4; P = 2 }
    5  = @{ B = 4; I = 1; P = 6 }
    6  = @{ I = 2; P = 2; AD = 1 }
    7  = @{ B = 3 }
    8  = @{ C = 5; P = 9; AK = 2 }
    11 = @{ AY = 1 }
    19 = @{ AD = 1 }
    20 = @{ AY = 1 }
    23 = @{ AK = 1 }
    25 = @{ I = 1 }
    26 = @{ P = 1 }
    29 = @{ AD = 2; AY = 1 }
    31 = @{ AK = 2 }
    35 = @{ I = 1 }
    41 = @{ AK = 1 }
    47 = @{ I = 1 }
    53 = @{ P = 6 }
    57 = @{ B = 2 }
    89 = @{ AR = 1 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
